# IPA New Script for Change Password in Account settings
# Adds a new test-case row (row 11) to the "Test Cases" sheet, including a
# hyperlink to the related JIRA issue, mirroring the pattern used by the
# existing rows (e.g. row 10 -> OPQA-4221).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Register the hyperlink for B11 first (this allocates rId3 in order) ---
$ws.Hyperlinks.Add($ws.Range("B11"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-4223", "", "", "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-4223")

# Adding the hyperlink forces Excel's built-in "Hyperlink" cell style onto
# B11 and registers that named style in the workbook. Remove the named style
# again so it does not linger in the saved styles table.
$wb.Styles.Item("Hyperlink").Delete()

# --- Copy formatting (borders/fill/font/alignment only - no values) from
#     existing cells that already carry the styles we need for row 11,
#     re-using the workbook's existing style entries instead of creating
#     new ones. ---
$ws.Range("A3").Copy()
$ws.Range("A11").PasteSpecial(-4122)

$ws.Range("C8").Copy()
$ws.Range("B11").PasteSpecial(-4122)

$ws.Range("B5").Copy()
$ws.Range("C11").PasteSpecial(-4122)

$ws.Range("A7").Copy()
$ws.Range("D11").PasteSpecial(-4122)

$ws.Range("E3").Copy()
$ws.Range("E11").PasteSpecial(-4122)

# --- Now fill in the real cell values. Order matches the order new shared
#     strings should be appended (B11, C11, A11, D11 - D11 re-uses the
#     existing "Y" string). ---
$ws.Range("B11").Value = "OPQA-4223 || OPQA-4224"
$ws.Range("C11").Value = "Verify that error message "" New password should not match current password"" should be displayed when user enters the current password in change password field.|| Verify that error message""New password should not match previous 4 passwords"" should be displayed when user enters password in change password field which is matching with the previous 4 passwords."
$ws.Range("A11").Value = "IPA0005"
$ws.Range("D11").Value = "Y"

# --- Row height to match the other wrapped-text rows (90pt) ---
$ws.Rows.Item(11).RowHeight = 90

# --- Update the visible selection to B11, matching the new active cell ---
$ws.Range("B11").Select()
